$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "Datos actualizados" timestamp message in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 14:05"

# 2. Row 28 (Suecia) - updated counts
$ws.Range("B28").Value = 35088
$ws.Range("C28").Value = 648
$ws.Range("E28").Value = 25897
$ws.Range("G28").Value = 95
$ws.Range("H28").Value = 4220

# 3. Row 41 (Rumania) - updated counts
$ws.Range("E41").Value = 5209
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 1223

# 4. Rows 49-50: Dinamarca moves above Panama with new Dinamarca data,
#    Panama keeps its old data but shifts to row 50
$ws.Range("A49").Value = "Dinamarca"
$ws.Range("B49").Value = 11480
$ws.Range("C49").Value = 52
$ws.Range("D49").Value = 10106
$ws.Range("E49").Value = 809
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 565

$ws.Range("A50").Value = "Panama"
$ws.Range("B50").Value = 11447
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 6379
$ws.Range("E50").Value = 4755
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 313

# 5. Row 76 (Uzbekistan) - updated counts
$ws.Range("B76").Value = 3355
$ws.Range("C76").Value = 65
$ws.Range("E76").Value = 682

# 6. Row 178 - updated counts
$ws.Range("D178").Value = 60
$ws.Range("E178").Value = 2

# 7. Row 195 (Laos) - updated counts
$ws.Range("D195").Value = 16
$ws.Range("E195").Value = 3
